$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was inserted before the existing row 222 ("Hortaliza,
# Femacal de La Calera - Zanahoria"), shifting all subsequent rows (222-316) down
# by one, growing the used range from A1:R316 to A1:R317.
$ws.Rows("222:222").Insert()

$ws.Range("A222").Value = 3
$ws.Range("B222").Value = "Femacal de La Calera"
$ws.Range("C222").Value = "Coquimbo"
$ws.Range("D222").Value = 44636
$ws.Range("E222").Value = 5
$ws.Range("F222").Value = 100114013
$ws.Range("G222").Value = "Zanahoria"
$ws.Range("H222").Value = "Sin especificar"
$ws.Range("I222").Value = "Primera"
$ws.Range("J222").Value = 370
$ws.Range("K222").Value = 7500
$ws.Range("L222").Value = 8000
$ws.Range("M222").Value = 7757
$ws.Range("N222").Value = "`$/saco 20 kilos"
$ws.Range("O222").Value = "Chillán"
$ws.Range("P222").Value = 388
$ws.Range("Q222").Value = 20
$ws.Range("R222").Value = "Hortaliza"
